# Change ILED to 0
# Add four new rows (8-11) mirroring rows 4-7 but for the NCAP_ILED
# attribute set to 0 instead of NCAP_START set to 2027.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$plants = @(
    "P-TH-CCGT-GAS-CCS04-Cork1",
    "P-TH-CCGT-GAS-CCS04-Cork2",
    "P-TH-CCGT-GAS-CCS04-Dublin1",
    "P-TH-CCGT-GAS-CCS04-Dublin2"
)

for ($i = 0; $i -lt 4; $i++) {
    $row = 8 + $i
    $ws.Cells.Item($row, 4).Value = "NCAP_ILED"
    $ws.Cells.Item($row, 8).Value = 0
    $ws.Cells.Item($row, 10).Value = $plants[$i]
    $ws.Cells.Item($row, 12).Value = "PWRGAS"
}

$ws.Range("D8:I11").Style = $ws.Range("D6:I6").Style
$ws.Range("K8:K11").Style = $ws.Range("K6:K6").Style
$ws.Range("J8:J11").Style = $ws.Range("J6:J6").Style
$ws.Range("L8:L11").Style = $ws.Range("L6:L6").Style

$ws.Range("M16").Select()
